# cn-#18 ignore columns after one empty cell in table header
$wb = $excel.ActiveWorkbook

$wsAdd = $wb.Worksheets.Item("AddTwoNumbers")
$wsSum = $wb.Worksheets.Item("SumProduct")

# Update the shared string text for the "Ignored after ..." header cell.
$wsSum.Range("F4").Value = "Ignored after an empty cells on a table"

# Move that cell one column to the left: F4 -> E4 (one empty cell gap instead of two).
$wsSum.Range("F4").Cut($wsSum.Range("E4"))

# Update the selection / active-cell on the SumProduct sheet.
$wsSum.Range("E5").Select()

# Make the AddTwoNumbers sheet the active (selected) tab instead of SumProduct.
$wsAdd.Activate()
